$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Price" (column D) and "Volume(1h)" (column E) figures for the
# cryptos list. In the source workbook every one of these cells is stored
# as plain text (inline strings such as "326.55" or "  -0.39%  "), so the
# replacement values must also land as text rather than being silently
# reinterpreted by Excel as numbers (which would both change the cell type
# and risk floating point rounding, e.g. "0.06970" -> 6.97E-02).
#
# For any replacement value that looks like a genuine number, we briefly
# switch the cell to a text NumberFormat before assigning it, then restore
# the cell's style to "Normal" so no stray formatting/style is left behind
# on the cell (matching the original, unstyled cells exactly).
$updates = @(
    @{ Cell = 'D2'; Value = '29.227.62' },
    @{ Cell = 'E2'; Value = '  +0.23%  ' },
    @{ Cell = 'D3'; Value = '1.900.40' },
    @{ Cell = 'E3'; Value = '  -0.44%  ' },
    @{ Cell = 'E4'; Value = '  -0.22%  ' },
    @{ Cell = 'D5'; Value = '326.55' },
    @{ Cell = 'E5'; Value = '  -0.39%  ' },
    @{ Cell = 'E6'; Value = '  -0.31%  ' },
    @{ Cell = 'D7'; Value = '0.4653' },
    @{ Cell = 'E7'; Value = '  +0.12%  ' },
    @{ Cell = 'D8'; Value = '0.3919' },
    @{ Cell = 'E8'; Value = '  -0.16%  ' },
    @{ Cell = 'D9'; Value = '0.07885' },
    @{ Cell = 'E9'; Value = '  -0.95%  ' },
    @{ Cell = 'D10'; Value = '0.9882' },
    @{ Cell = 'E10'; Value = '  -1.58%  ' },
    @{ Cell = 'D11'; Value = '21.97' },
    @{ Cell = 'E11'; Value = '  -1.36%  ' },
    @{ Cell = 'D12'; Value = '1.917.20' },
    @{ Cell = 'E12'; Value = '  +0.08%  ' },
    @{ Cell = 'D13'; Value = '7.077' },
    @{ Cell = 'E13'; Value = '  -0.77%  ' },
    @{ Cell = 'D14'; Value = '5.750' },
    @{ Cell = 'E14'; Value = '  -0.71%  ' },
    @{ Cell = 'D15'; Value = '0.06970' },
    @{ Cell = 'E15'; Value = '  -0.17%  ' },
    @{ Cell = 'D16'; Value = '88.35' },
    @{ Cell = 'E16'; Value = '  -0.32%  ' },
    @{ Cell = 'D17'; Value = '1.003' },
    @{ Cell = 'E17'; Value = '  -0.20%  ' },
    @{ Cell = 'D18'; Value = '0.000009981' },
    @{ Cell = 'E18'; Value = '  -1.26%  ' },
    @{ Cell = 'D19'; Value = '17.08' },
    @{ Cell = 'E19'; Value = '  -1.02%  ' },
    @{ Cell = 'E20'; Value = '  -0.30%  ' },
    @{ Cell = 'D21'; Value = '29.249.97' },
    @{ Cell = 'E21'; Value = '  +0.21%  ' },
    @{ Cell = 'D22'; Value = '5.317' },
    @{ Cell = 'E22'; Value = '  -1.01%  ' },
    @{ Cell = 'E23'; Value = '  +0.10%  ' },
    @{ Cell = 'D24'; Value = '2.095' },
    @{ Cell = 'E24'; Value = '  +1.83%  ' },
    @{ Cell = 'D25'; Value = '156.29' },
    @{ Cell = 'E25'; Value = '  -0.04%  ' },
    @{ Cell = 'D26'; Value = '19.45' },
    @{ Cell = 'E26'; Value = '  -0.42%  ' },
    @{ Cell = 'D27'; Value = '5.983' },
    @{ Cell = 'E27'; Value = '  +2.37%  ' },
    @{ Cell = 'D28'; Value = '118.62' },
    @{ Cell = 'E28'; Value = '  -0.82%  ' },
    @{ Cell = 'D29'; Value = '1.909' },
    @{ Cell = 'E29'; Value = '  -4.65%  ' },
    @{ Cell = 'D30'; Value = '0.09356' },
    @{ Cell = 'E30'; Value = '  -0.43%  ' },
    @{ Cell = 'D31'; Value = '0.9061' },
    @{ Cell = 'E31'; Value = '  -1.81%  ' },
    @{ Cell = 'D32'; Value = '5.285' },
    @{ Cell = 'E32'; Value = '  -1.52%  ' },
    @{ Cell = 'D33'; Value = '1.325' },
    @{ Cell = 'E33'; Value = '  -1.39%  ' },
    @{ Cell = 'D34'; Value = '3.216' },
    @{ Cell = 'E34'; Value = '  -1.84%  ' },
    @{ Cell = 'D35'; Value = '1.183' },
    @{ Cell = 'E35'; Value = '  +1.97%  ' },
    @{ Cell = 'D36'; Value = '0.05782' },
    @{ Cell = 'E36'; Value = '  -1.02%  ' },
    @{ Cell = 'D37'; Value = '0.02089' },
    @{ Cell = 'E37'; Value = '  -0.32%  ' },
    @{ Cell = 'D38'; Value = '1.001' },
    @{ Cell = 'E38'; Value = '  -0.31%  ' },
    @{ Cell = 'D39'; Value = '7.748' },
    @{ Cell = 'E39'; Value = '  -3.38%  ' },
    @{ Cell = 'D40'; Value = '0.5711' },
    @{ Cell = 'E40'; Value = '  -0.73%  ' },
    @{ Cell = 'E41'; Value = '  -1.18%  ' },
    @{ Cell = 'D42'; Value = '9.754' },
    @{ Cell = 'E42'; Value = '  -2.44%  ' },
    @{ Cell = 'D43'; Value = '11.95' },
    @{ Cell = 'E43'; Value = '  -0.39%  ' },
    @{ Cell = 'D44'; Value = '0.5350' },
    @{ Cell = 'E44'; Value = '  -1.38%  ' },
    @{ Cell = 'D45'; Value = '2.196' },
    @{ Cell = 'E45'; Value = '  -0.99%  ' },
    @{ Cell = 'D46'; Value = '0.07046' },
    @{ Cell = 'E46'; Value = '  -0.68%  ' },
    @{ Cell = 'D47'; Value = '1.857' },
    @{ Cell = 'E47'; Value = '  -1.44%  ' },
    @{ Cell = 'D48'; Value = '2.577' },
    @{ Cell = 'E48'; Value = '  -0.50%  ' },
    @{ Cell = 'D49'; Value = '113.25' },
    @{ Cell = 'E49'; Value = '  +0.81%  ' },
    @{ Cell = 'D50'; Value = '1.058' },
    @{ Cell = 'E50'; Value = '  -1.34%  ' },
    @{ Cell = 'D51'; Value = '71.19' },
    @{ Cell = 'E51'; Value = '  -0.53%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)

    if ($u.Value -match '^[+-]?[0-9]*\.?[0-9]+$') {
        $cell.NumberFormat = '@'
        $cell.Value = $u.Value
        $cell.Style = 'Normal'
    } else {
        $cell.Value = $u.Value
    }
}
